$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had two "normal" COD rows (3 and 5) and two
# out-of-range / flagged COD rows (4 and 6, highlighted with a red fill).
# All four are being dropped from the report, which shifts the trailing
# "Average = " (was row 8) and "COD COUNT = " (was row 10) rows up to
# rows 4 and 6 respectively.
$ws.Range("A3:A6").EntireRow.Delete()

# With the two flagged rows gone, recompute the summary figures from the
# remaining data (only row 2 is left): the average of {6.62} is 6.62, and
# the COD count (entries whose id contains "COD") is still 0.
$ws.Range("B4").Value = 6.62
$ws.Range("B6").Value = 0
